{"js": "const replacements = [\n  [\"981\u00f75=196, 1\", \"508\u00f73=169, 1\"],\n  [\"340\u00f72=170, 0\", \"789\u00f76=131, 3\"],\n  [\"289\u00f74=72, 1\", \"762\u00f73=254, 0\"],\n  [\"339\u00f72=169, 1\", \"714\u00f72=357, 0\"],\n  [\"942\u00f76=157, 0\", \"612\u00f72=306, 0\"],\n  [\"503\u00f75=100, 3\", \"594\u00f73=198, 0\"],\n  [\"972\u00f72=486, 0\", \"526\u00f72=263, 0\"],\n  [\"672\u00f79=74, 6\", \"699\u00f75=139, 4\"],\n  [\"862\u00f77=123, 1\", \"818\u00f76=136, 2\"],\n  [\"793\u00f79=88, 1\", \"641\u00f79=71, 2\"],\n  [\"717\u00f76=119, 3\", \"317\u00f72=158, 1\"],\n  [\"113\u00f72=56, 1\", \"169\u00f77=24, 1\"],\n  [\"633\u00f73=211, 0\", \"816\u00f79=90, 6\"],\n  [\"936\u00f74=234, 0\", \"605\u00f76=100, 5\"],\n  [\"418\u00f79=46, 4\", \"665\u00f79=73, 8\"],\n  [\"768\u00f79=85, 3\", \"408\u00f72=204, 0\"],\n  [\"391\u00f73=130, 1\", \"915\u00f74=228, 3\"],\n  [\"277\u00f75=55, 2\", \"146\u00f73=48, 2\"],\n  [\"848\u00f73=282, 2\", \"563\u00f74=140, 3\"],\n  [\"889\u00f75=177, 4\", \"133\u00f72=66, 1\"],\n  [\"574\u00f78=71, 6\", \"659\u00f78=82, 3\"],\n  [\"668\u00f79=74, 2\", \"751\u00f73=250, 1\"],\n  [\"552\u00f76=92, 0\", \"738\u00f76=123, 0\"],\n  [\"816\u00f74=204, 0\", \"791\u00f75=158, 1\"],\n  [\"229\u00f74=57, 1\", \"146\u00f75=29, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"981\u00f75=196, 1\", \"508\u00f73=169, 1\"),\n  @(\"340\u00f72=170, 0\", \"789\u00f76=131, 3\"),\n  @(\"289\u00f74=72, 1\", \"762\u00f73=254, 0\"),\n  @(\"339\u00f72=169, 1\", \"714\u00f72=357, 0\"),\n  @(\"942\u00f76=157, 0\", \"612\u00f72=306, 0\"),\n  @(\"503\u00f75=100, 3\", \"594\u00f73=198, 0\"),\n  @(\"972\u00f72=486, 0\", \"526\u00f72=263, 0\"),\n  @(\"672\u00f79=74, 6\", \"699\u00f75=139, 4\"),\n  @(\"862\u00f77=123, 1\", \"818\u00f76=136, 2\"),\n  @(\"793\u00f79=88, 1\", \"641\u00f79=71, 2\"),\n  @(\"717\u00f76=119, 3\", \"317\u00f72=158, 1\"),\n  @(\"113\u00f72=56, 1\", \"169\u00f77=24, 1\"),\n  @(\"633\u00f73=211, 0\", \"816\u00f79=90, 6\"),\n  @(\"936\u00f74=234, 0\", \"605\u00f76=100, 5\"),\n  @(\"418\u00f79=46, 4\", \"665\u00f79=73, 8\"),\n  @(\"768\u00f79=85, 3\", \"408\u00f72=204, 0\"),\n  @(\"391\u00f73=130, 1\", \"915\u00f74=228, 3\"),\n  @(\"277\u00f75=55, 2\", \"146\u00f73=48, 2\"),\n  @(\"848\u00f73=282, 2\", \"563\u00f74=140, 3\"),\n  @(\"889\u00f75=177, 4\", \"133\u00f72=66, 1\"),\n  @(\"574\u00f78=71, 6\", \"659\u00f78=82, 3\"),\n  @(\"668\u00f79=74, 2\", \"751\u00f73=250, 1\"),\n  @(\"552\u00f76=92, 0\", \"738\u00f76=123, 0\"),\n  @(\"816\u00f74=204, 0\", \"791\u00f75=158, 1\"),\n  @(\"229\u00f74=57, 1\", \"146\u00f75=29, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $old\n    $rng.Find.Replacement.Text = $new\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Execute([ref]$old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2)\n}\n"}
